$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table by one more year column (O), mirroring the format of
# the existing last column (N) for the header band (row 3), the year
# header (row 4) and the data rows (5-10).
$ws.Range("N3:N10").Copy() | Out-Null
$ws.Range("O3:O10").PasteSpecial(-4122) | Out-Null

# New year header
$ws.Range("O4").Value = 2021

# New data values for the 2021 column
$ws.Range("O6").Value = 1860
$ws.Range("O7").Value = 1
$ws.Range("O8").Value = 510
$ws.Range("O9").Value = 178
$ws.Range("O10").Value = 821

# Update the recorded selection to match the post-edit state
$ws.Range("P9").Select() | Out-Null
